$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# 1. Copy the formatting of the existing placeholder row (row 6) down onto
#    the rows that will hold the six new application entries (rows 6-10 will
#    become the new app rows, rows 11-12 will re-host the placeholder entry).
#    Using PasteSpecial(xlPasteFormats) re-uses the existing cell styles
#    (s="4"/"3"/"2") instead of minting brand-new style entries.
# ---------------------------------------------------------------------------

$ws.Range("A6:D6").Copy()
$ws.Range("A7:D12").PasteSpecial(-4122)

# Column C on the new app rows (7-10) is a real hyperlink cell, so borrow the
# hyperlink-flavoured style (s="7") that is already used by C3:C5, instead of
# the plain style used by the placeholder's C6.
$ws.Range("C5").Copy()
$ws.Range("C7:C10").PasteSpecial(-4122)

# ---------------------------------------------------------------------------
# 2. Fill in the values for the six new rows of data.
# ---------------------------------------------------------------------------

$ws.Range("A6").Value = "History Skeleton Generator"
$ws.Range("B6").Value = "This application is writen using Jupyter Notebook. Inside the notebook are the instructions for creating a Timeline workbook for the program to use to create an output workbook containing rough timelines for nations and factions in a fictional universe."
$ws.Range("C6").Value = "https://github.com/valenpendragon/history-skeleton-generator"
$ws.Range("D6").Value = "5.png"

$ws.Range("A7").Value = "Web Weather API"
$ws.Range("B7").Value = "This application is a Flask website that provides climate data to the user using a URL API. Instructions for using the API are on the site's homepage."
$ws.Range("C7").Value = "https://github.com/valenpendragon/web-weather-api"
$ws.Range("D7").Value = "5.png"

$ws.Range("A8").Value = "NASA's Astronony Web Page"
$ws.Range("B8").Value = "This is a simple Streamlit App that displays NASA's Astronomy Picture of the Day, including the Copyright information when it exists."
$ws.Range("C8").Value = "https://github.com/valenpendragon/apod-web-page"
$ws.Range("D8").Value = "5.png"

$ws.Range("A9").Value = "News API Email"
$ws.Range("B9").Value = "This application gathers daily news on a specified topic and emails the article links to the User. It serves more to demonstrate how to add this functionality to a web site."
$ws.Range("C9").Value = "https://github.com/valenpendragon/news-api-email"
$ws.Range("D9").Value = "5.png"

$ws.Range("A10").Value = "Spreadsheet to PDF Invoice generator"
$ws.Range("B10").Value = "This is a backend application that takes spreadsheet invoices generated from another part of an application and turns them into PDF Invoices that can be emailed to customers."
$ws.Range("C10").Value = "https://github.com/valenpendragon/PDF-invoices"
$ws.Range("D10").Value = "5.png"

# Rows 11 and 12 re-host the original placeholder entry that used to live on
# row 6 (it is duplicated onto two rows in the new layout).
$ws.Range("A11").Value = "Placeholder2"
$ws.Range("B11").Value = "This is a placeholder to ensure the portfolio app works."
$ws.Range("C11").Value = "placedholder"
$ws.Range("D11").Value = "5.png"

$ws.Range("A12").Value = "Placeholder2"
$ws.Range("B12").Value = "This is a placeholder to ensure the portfolio app works."
$ws.Range("C12").Value = "placedholder"
$ws.Range("D12").Value = "5.png"

# ---------------------------------------------------------------------------
# 3. Row heights for the wrapped-text rows.
# ---------------------------------------------------------------------------

$ws.Rows(6).RowHeight = 60
$ws.Rows(7).RowHeight = 45
$ws.Rows(8).RowHeight = 45
$ws.Rows(9).RowHeight = 45
$ws.Rows(10).RowHeight = 45

# ---------------------------------------------------------------------------
# 4. Hyperlinks for the new URL cells.
# ---------------------------------------------------------------------------

$ws.Hyperlinks.Add($ws.Range("C6"), "https://github.com/valenpendragon/history-skeleton-generator")
$ws.Hyperlinks.Add($ws.Range("C7"), "https://github.com/valenpendragon/web-weather-api")
$ws.Hyperlinks.Add($ws.Range("C8"), "https://github.com/valenpendragon/apod-web-page")
$ws.Hyperlinks.Add($ws.Range("C9"), "https://github.com/valenpendragon/news-api-email")
$ws.Hyperlinks.Add($ws.Range("C10"), "https://github.com/valenpendragon/PDF-invoices")

# ---------------------------------------------------------------------------
# 5. Leave the active selection on C8, matching the author's last-saved view.
# ---------------------------------------------------------------------------

$ws.Range("C8").Select()
